# Updated cryptos list on Tue Apr 11 10:43:42 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = "30.118.23"
$ws.Range("E2").Value = "  +5.77%  "

# --- Row 3 ---
$ws.Range("D3").Value = "1.919.68"
$ws.Range("E3").Value = "  +2.74%  "

# --- Row 4 ---
$ws.Range("E4").Value = "  -0.68%  "

# --- Row 5 ---
$ws.Range("D5").Value = "'330.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.74%  "

# --- Row 6 ---
$ws.Range("D6").Value = "'0.9999"
$ws.Range("D6").Style = "Normal"

# --- Row 7 ---
$ws.Range("D7").Value = "'0.5214"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.93%  "

# --- Row 8 ---
$ws.Range("D8").Value = "'0.4095"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.99%  "

# --- Row 9 ---
$ws.Range("D9").Value = "'0.08520"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.92%  "

# --- Row 10 ---
$ws.Range("D10").Value = "'1.128"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.08%  "

# --- Row 11 ---
$ws.Range("D11").Value = "'42.84"
$ws.Range("D11").Style = "Normal"

# --- Row 12 ---
$ws.Range("E12").Value = "  +9.85%  "

# --- Row 13 ---
$ws.Range("D13").Value = "'6.435"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.63%  "

# --- Row 14 ---
$ws.Range("D14").Value = "1.909.23"
$ws.Range("E14").Value = "  +1.58%  "

# --- Row 15 ---
$ws.Range("D15").Value = "'7.412"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.97%  "

# --- Row 16 ---
$ws.Range("E16").Value = "  -0.74%  "

# --- Row 17 ---
$ws.Range("D17").Value = "'95.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.94%  "

# --- Row 18 ---
$ws.Range("E18").Value = "  +1.11%  "

# --- Row 19 ---
$ws.Range("D19").Value = "'0.06690"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.65%  "

# --- Row 20 ---
$ws.Range("D20").Value = "'18.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.76%  "

# --- Row 21 ---
$ws.Range("D21").Value = "'0.9997"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.62%  "

# --- Row 22 ---
$ws.Range("D22").Value = "'6.014"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.59%  "

# --- Row 23 ---
$ws.Range("D23").Value = "30.125.24"
$ws.Range("E23").Value = "  +5.68%  "

# --- Row 24 ---
$ws.Range("D24").Value = "'11.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.21%  "

# --- Row 25 ---
$ws.Range("D25").Value = "'2.213"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.20%  "

# --- Row 26 ---
$ws.Range("D26").Value = "2.136.53"
$ws.Range("E26").Value = "  +2.09%  "

# --- Row 27 ---
$ws.Range("D27").Value = "'21.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.52%  "

# --- Row 28 ---
$ws.Range("D28").Value = "'159.92"
$ws.Range("D28").Style = "Normal"

# --- Row 29 ---
$ws.Range("D29").Value = "'2.448"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.78%  "

# --- Row 30 ---
$ws.Range("D30").Value = "'129.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.84%  "

# --- Row 31 ---
$ws.Range("E31").Value = "  +3.95%  "

# --- Row 32 ---
$ws.Range("E32").Value = "  +1.72%  "

# --- Row 33 ---
$ws.Range("D33").Value = "'6.068"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.86%  "

# --- Row 34 ---
$ws.Range("D34").Value = "'3.634"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.43%  "

# --- Row 35 ---
$ws.Range("D35").Value = "'0.02489"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.36%  "

# --- Row 36 ---
$ws.Range("D36").Value = "'0.06631"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.51%  "

# --- Row 37 ---
$ws.Range("E37").Value = "  +2.42%  "

# --- Row 38 ---
$ws.Range("E38").Value = "  +4.99%  "

# --- Row 39 ---
$ws.Range("D39").Value = "'5.199"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.52%  "

# --- Row 40 ---
$ws.Range("D40").Value = "'8.903"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.23%  "

# --- Row 41 ---
$ws.Range("D41").Value = "'0.6560"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.05%  "

# --- Row 42 ---
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.250"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.55%  "

# --- Row 43 ---
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "'11.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.09%  "

# --- Row 44 ---
$ws.Range("D44").Value = "'0.6164"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.64%  "

# --- Row 45 ---
$ws.Range("D45").Value = "'13.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.57%  "

# --- Row 46 ---
$ws.Range("D46").Value = "'3.772"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.42%  "

# --- Row 47 ---
$ws.Range("D47").Value = "'2.086"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.03%  "

# --- Row 48 ---
$ws.Range("D48").Value = "'1.251"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.03%  "

# --- Row 49 ---
$ws.Range("D49").Value = "'124.70"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.90%  "

# --- Row 50 ---
$ws.Range("D50").Value = "'1.177"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +10.24%  "

# --- Row 51 ---
$ws.Range("D51").Value = "'79.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.15%  "

